# Updates crypto price/volume figures in Sheet1 (columns D and E, rows 2-51)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "45.569.27"
$ws.Range("E2").Value = "  +6.84%  "
$ws.Range("D3").Value = "2.397.47"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'114.88"
$ws.Range("E5").Value = "  +9.93%  "
$ws.Range("D6").Value = "'319.05"
$ws.Range("E6").Value = "  +2.82%  "
$ws.Range("E7").Value = "  +2.91%  "
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("D9").Value = "'0.628"
$ws.Range("E9").Value = "  +3.84%  "
$ws.Range("D10").Value = "'42.58"
$ws.Range("E10").Value = "  +7.84%  "
$ws.Range("D11").Value = "'0.0931"
$ws.Range("E11").Value = "  +3.04%  "
$ws.Range("D12").Value = "'8.69"
$ws.Range("E12").Value = "  +5.27%  "
$ws.Range("D13").Value = "'0.110"
$ws.Range("E13").Value = "  +2.83%  "
$ws.Range("E14").Value = "  +3.08%  "
$ws.Range("D15").Value = "'16.00"
$ws.Range("E15").Value = "  +4.23%  "
$ws.Range("D16").Value = "2.760.66"
$ws.Range("E16").Value = "  -0.82%  "
$ws.Range("D17").Value = "2.401.27"
$ws.Range("D18").Value = "45.579.18"
$ws.Range("E18").Value = "  +6.41%  "
$ws.Range("D19").Value = "'7.51"
$ws.Range("E19").Value = "  +2.61%  "
$ws.Range("D20").Value = "'0.0000108"
$ws.Range("E20").Value = "  +3.39%  "
$ws.Range("D21").Value = "'13.57"
$ws.Range("E21").Value = "  +1.03%  "
$ws.Range("D22").Value = "'74.95"
$ws.Range("E22").Value = "  +2.23%  "
$ws.Range("D23").Value = "'3.59"
$ws.Range("E23").Value = "  +4.26%  "
$ws.Range("D24").Value = "'265.04"
$ws.Range("E24").Value = "  -1.12%  "
$ws.Range("D25").Value = "'2.36"
$ws.Range("E25").Value = "  +6.47%  "
$ws.Range("E26").Value = "  -0.80%  "
$ws.Range("D27").Value = "'7.75"
$ws.Range("E27").Value = "  +5.92%  "
$ws.Range("D28").Value = "'11.40"
$ws.Range("E28").Value = "  +4.76%  "
$ws.Range("E29").Value = "  +2.57%  "
$ws.Range("D30").Value = "'39.90"
$ws.Range("E30").Value = "  +10.18%  "
$ws.Range("E31").Value = "  +15.71%  "
$ws.Range("D32").Value = "'22.85"
$ws.Range("E32").Value = "  +2.57%  "
$ws.Range("D33").Value = "'173.11"
$ws.Range("E33").Value = "  +4.84%  "
$ws.Range("E34").Value = "  +11.64%  "
$ws.Range("E35").Value = "  +1.89%  "
$ws.Range("D36").Value = "'5.02"
$ws.Range("E36").Value = "  +10.49%  "
$ws.Range("D37").Value = "'0.118"
$ws.Range("E37").Value = "  +6.78%  "
$ws.Range("D38").Value = "'4.18"
$ws.Range("E38").Value = "  +15.38%  "
$ws.Range("E39").Value = "  +10.02%  "
$ws.Range("D40").Value = "'0.0365"
$ws.Range("E40").Value = "  +5.15%  "
$ws.Range("D41").Value = "'1.79"
$ws.Range("E41").Value = "  +13.12%  "
$ws.Range("E42").Value = "  +7.22%  "
$ws.Range("D43").Value = "'13.69"
$ws.Range("E43").Value = "  +11.67%  "
$ws.Range("D44").Value = "'100.36"
$ws.Range("E44").Value = "  -8.58%  "
$ws.Range("D45").Value = "'72.02"
$ws.Range("E45").Value = "  +1.26%  "
$ws.Range("D46").Value = "'87.84"
$ws.Range("E46").Value = "  +13.38%  "
$ws.Range("D48").Value = "'5.86"
$ws.Range("E48").Value = "  +13.98%  "
$ws.Range("D49").Value = "'116.25"
$ws.Range("E49").Value = "  +4.79%  "
$ws.Range("E50").Value = "  +9.85%  "
$ws.Range("D51").Value = "'1.58"
$ws.Range("E51").Value = "  +10.95%  "
